$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = 111
$ws.Range("B6").Value = 222
$ws.Range("C6").Value = 333
$ws.Range("A7").Value = 888
$ws.Range("B7").Value = 9999
$ws.Range("B7").Select()
